$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on affected columns so values stay text (matching original inlineStr type)
$colD = $ws.Range("D2:D51")
$colE = $ws.Range("E2:E51")
$colG = $ws.Range("G2:G51")
$colD.NumberFormat = "@"
$colE.NumberFormat = "@"
$colG.NumberFormat = "@"

$ws.Range("D2").Value = "304.94"
$ws.Range("E2").Value = "1.19%"
$ws.Range("G2").Value = "18"
$ws.Range("D3").Value = "35.86"
$ws.Range("E3").Value = "1.13%"
$ws.Range("G3").Value = "18"
$ws.Range("D4").Value = "5.073"
$ws.Range("E4").Value = "0.07%"
$ws.Range("G4").Value = "18"
$ws.Range("D5").Value = "0.08055"
$ws.Range("E5").Value = "0.94%"
$ws.Range("G5").Value = "18"
$ws.Range("D6").Value = "1.934"
$ws.Range("E6").Value = "2.76%"
$ws.Range("G6").Value = "18"
$ws.Range("D7").Value = "4.167"
$ws.Range("E7").Value = "2.99%"
$ws.Range("G7").Value = "18"
$ws.Range("D8").Value = "7.847"
$ws.Range("E8").Value = "0.76%"
$ws.Range("G8").Value = "18"
$ws.Range("D9").Value = "0.9317"
$ws.Range("E9").Value = "0.37%"
$ws.Range("G9").Value = "18"
$ws.Range("D10").Value = "0.1263"
$ws.Range("E10").Value = "-12.51%"
$ws.Range("G10").Value = "18"
$ws.Range("D11").Value = "0.1913"
$ws.Range("E11").Value = "-0.73%"
$ws.Range("G11").Value = "18"
$ws.Range("D12").Value = "0.09188"
$ws.Range("E12").Value = "1.00%"
$ws.Range("G12").Value = "18"
$ws.Range("D13").Value = "0.03490"
$ws.Range("E13").Value = "-0.29%"
$ws.Range("G13").Value = "18"
$ws.Range("D14").Value = "0.09905"
$ws.Range("E14").Value = "0.63%"
$ws.Range("G14").Value = "18"
$ws.Range("D15").Value = "0.001434"
$ws.Range("E15").Value = "2.51%"
$ws.Range("G15").Value = "18"
$ws.Range("D16").Value = "0.006639"
$ws.Range("E16").Value = "13.08%"
$ws.Range("G16").Value = "18"
$ws.Range("D17").Value = "3.613"
$ws.Range("E17").Value = "2.29%"
$ws.Range("G17").Value = "18"
$ws.Range("D18").Value = "3.168"
$ws.Range("E18").Value = "7.03%"
$ws.Range("G18").Value = "18"
$ws.Range("D19").Value = "0.3422"
$ws.Range("E19").Value = "-0.06%"
$ws.Range("G19").Value = "18"
$ws.Range("E20").Value = "2.62%"
$ws.Range("G20").Value = "18"
$ws.Range("D21").Value = "5.189"
$ws.Range("E21").Value = "2.62%"
$ws.Range("G21").Value = "18"
$ws.Range("D22").Value = "0.2537"
$ws.Range("E22").Value = "5.78%"
$ws.Range("G22").Value = "18"
$ws.Range("D23").Value = "0.04407"
$ws.Range("E23").Value = "-1.93%"
$ws.Range("G23").Value = "18"
$ws.Range("D24").Value = "0.001238"
$ws.Range("E24").Value = "2.17%"
$ws.Range("G24").Value = "18"
$ws.Range("D25").Value = "0.004714"
$ws.Range("E25").Value = "-1.04%"
$ws.Range("G25").Value = "18"
$ws.Range("D26").Value = "0.0001304"
$ws.Range("E26").Value = "6.02%"
$ws.Range("G26").Value = "18"
$ws.Range("D27").Value = "0.0003137"
$ws.Range("E27").Value = "3.75%"
$ws.Range("G27").Value = "18"
$ws.Range("G28").Value = "18"
$ws.Range("G29").Value = "18"
$ws.Range("G30").Value = "18"
$ws.Range("G31").Value = "18"
$ws.Range("G32").Value = "18"
$ws.Range("G33").Value = "18"
$ws.Range("G34").Value = "18"
$ws.Range("G35").Value = "18"
$ws.Range("G36").Value = "18"
$ws.Range("G37").Value = "18"
$ws.Range("G38").Value = "18"
$ws.Range("D39").Value = "0.01977"
$ws.Range("E39").Value = "7.88%"
$ws.Range("G39").Value = "18"
$ws.Range("D40").Value = "0.05144"
$ws.Range("E40").Value = "8.22%"
$ws.Range("G40").Value = "18"
$ws.Range("D41").Value = "0.007604"
$ws.Range("E41").Value = "3.13%"
$ws.Range("G41").Value = "18"
$ws.Range("D42").Value = "0.01005"
$ws.Range("E42").Value = "-4.46%"
$ws.Range("G42").Value = "18"
$ws.Range("D43").Value = "0.1367"
$ws.Range("E43").Value = "2.86%"
$ws.Range("G43").Value = "18"
$ws.Range("D44").Value = "0.002106"
$ws.Range("E44").Value = "-0.16%"
$ws.Range("G44").Value = "18"
$ws.Range("E45").Value = "-2.13%"
$ws.Range("G45").Value = "18"
$ws.Range("D46").Value = "0.00006327"
$ws.Range("E46").Value = "1.58%"
$ws.Range("G46").Value = "18"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.29%"
$ws.Range("G47").Value = "18"
$ws.Range("G48").Value = "18"
$ws.Range("D49").Value = "0.001604"
$ws.Range("E49").Value = "-3.32%"
$ws.Range("G49").Value = "18"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").Value = "0.29%"
$ws.Range("G50").Value = "18"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").Value = "0.29%"
$ws.Range("G51").Value = "18"

# Restore default (Normal) style so no stray number-format styling is introduced
$colD.Style = "Normal"
$colE.Style = "Normal"
$colG.Style = "Normal"
